# "Snelle commit voor batterij op is!" - fill in manufacturer/model_name
# columns (parsed from the auction description text) and touch the
# remaining per-lot detail columns that the scraper emits but left blank
# for these rows (year, reference_number, material, case_number,
# diameter, movement_number, calibre, bracelet_strap, accessoires,
# signed) so they exist as blank cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$manufacturers = @{
    2 = "NATIONAL WATCH"
    3 = "GALLET"
    4 = "HELVETIA"
}

# Columns that should end up present-but-blank for each of these rows.
$blankCols = @("D", "E", "L", "M", "O", "P", "Q", "R", "S", "T")

foreach ($row in 2..4) {
    $name = $manufacturers[$row]

    $ws.Range("C$row").Value = $name
    $ws.Range("F$row").Value = $name

    foreach ($col in $blankCols) {
        $cell = $ws.Range("$col$row")
        # Assigning "" alone clears/removes the cell entirely in this
        # engine. A lone apostrophe is Excel's "force text" prefix and
        # is stripped from the stored value, so this yields a real,
        # present cell holding an empty TEXT string (matching the
        # scraper's empty inline/shared string cells) instead of no
        # cell at all.
        $cell.Value = "'"
        # The apostrophe prefix also flips on the "stored as text"
        # quote-prefix number format; re-apply the plain default style
        # so the cell stays style-neutral, same as the source cells.
        $cell.Style = "Normal"
    }
}
